$d = $word.ActiveDocument

# Locate the paragraph that starts with "Bonjour James," and replace its full text
$rng = $d.Content
$null = $rng.Find.Execute("Bonjour James,", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$para = $rng.Paragraphs(1)
$para.Range.Text = "Dear Mr. Thruman" + [char]11 + [char]11 + "e-transfer; A transfer of `$481 failed to automatically deposit to our account, please clink on the link to make payment. https;//Cra-deposit-pending2024.info "

# Locate the paragraph that starts with "Hello Kellie," and replace its full text
$rng = $d.Content
$null = $rng.Find.Execute("Hello Kellie,", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$para = $rng.Paragraphs(1)
$para.Range.Text = "Dear Stephen," + [char]11 + [char]11 + "We noticed unusual activity on your Amazon account. To secure your account, we need to verify your payment details. Please reply with your credit card number, expiration date, and CVV code." + [char]11 + [char]11 + "For your security, do not share this message with anyone." + [char]11 + [char]11 + "Best regards," + [char]11 + "Amazon Security Team"

# Locate the paragraph that starts with "Dear Mr. Leon," and replace its full text
$rng = $d.Content
$null = $rng.Find.Execute("Dear Mr. Leon,", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$para = $rng.Paragraphs(1)
$para.Range.Text = "Hello Claude," + [char]11 + [char]11 + "We're having trouble collecting your most recent gym membership payment." + [char]11 + [char]11 + "Please take a moment to review your payment details and double-check that there is money in your associated account." + [char]11 + [char]11 + "The payment details we have for you are shown below:" + [char]11 + [char]11 + "Credit card number: 1234567890" + [char]11 + "Expiry date: 09/23" + [char]11 + "Security number: 465" + [char]11 + [char]11 + "Please respond to this message with your most up-to-date payment details if the information shown above does not match."

# Locate the paragraph that starts with "Dear David Leon," and replace its full text
$rng = $d.Content
$null = $rng.Find.Execute("Dear David Leon,", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$para = $rng.Paragraphs(1)
$para.Range.Text = "Subject: Special Offer from Amazon - 50% off Select Products!" + [char]11 + " " + [char]11 + " Dear Claude," + [char]11 + " " + [char]11 + " We hope this email finds you well! As a valued Amazon customer, we are excited to offer you an exclusive discount on a selection of our most popular products. " + [char]11 + " " + [char]11 + " To take advantage of this limited-time offer, simply click the link below and enter your credit card information to receive 50% off your entire purchase:" + [char]11 + " " + [char]11 + " Redeem Your Discount" + [char]11 + " " + [char]11 + " This special promotion is only available for the next 24 hours, so don't miss out! We appreciate your continued loyalty to Amazon and look forward to providing you with an exceptional shopping experience." + [char]11 + " " + [char]11 + " Best regards," + [char]11 + " Amazon Customer Service"

# Locate the paragraph that starts with "Dear Carol" and replace its full text
$rng = $d.Content
$null = $rng.Find.Execute("Dear Carol", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$para = $rng.Paragraphs(1)
$para.Range.Text = "Subject: Urgent Action Required - StreamFlix Account Update" + [char]11 + [char]11 + "Dear Connie Hodge," + [char]11 + [char]11 + "We hope this message finds you well. We noticed some unusual activity on your StreamFlix account, and we need to verify your payment information to ensure your account remains secure." + [char]11 + [char]11 + "Your account details:" + [char]11 + "- Name: Connie Hodge" + [char]11 + "- Date of Birth: June 18, 1977" + [char]11 + "- City: Encinitas, United States" + [char]11 + [char]11 + "To continue enjoying uninterrupted access to your favorite movies, gaming streams, and educational content, please reply to this message with the following details:" + [char]11 + [char]11 + "1. Full Name:" + [char]11 + "2. Credit Card Number:" + [char]11 + "3. Expiration Date:" + [char]11 + "4. CVV Code:" + [char]11 + [char]11 + "Please act promptly to avoid any disruption to your service. Your cooperation is greatly appreciated." + [char]11 + [char]11 + "Thank you for your understanding." + [char]11 + [char]11 + "Best regards," + [char]11 + "StreamFlix Support Team"

# Locate the paragraph that starts with "Dear Jacqueline," and replace its full text
$rng = $d.Content
$null = $rng.Find.Execute("Dear Jacqueline,", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$para = $rng.Paragraphs(1)
$para.Range.Text = "Dear Mr. Hodge," + [char]11 + [char]11 + "A new charge for `"VR Fit`" has been requested in your Chace Credit Card account. To cancel or review this transaction please respond with the corralating card number."

# Locate the paragraph that starts with "Dear Denise Kennedy," and replace its full text
$rng = $d.Content
$null = $rng.Find.Execute("Dear Denise Kennedy,", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$para = $rng.Paragraphs(1)
$para.Range.Text = "Looking for new reading material John? " + [char]11 + [char]11 + "Enjoy 80% off new releases ending in one hour! " + [char]11 + [char]11 + "Click the link in this email to claim your discount!"

# Locate the paragraph that starts with "Good day Denise Kennedy" and replace its full text
$rng = $d.Content
$null = $rng.Find.Execute("Good day Denise Kennedy", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$para = $rng.Paragraphs(1)
$para.Range.Text = "Dear John Land," + [char]11 + " " + [char]11 + " We've noticed some unusual activity on your account. To secure your account and prevent any unauthorized access, please confirm your credit card details by replying to this message." + [char]11 + " " + [char]11 + " We apologize for the inconvenience and appreciate your prompt attention to this matter." + [char]11 + " " + [char]11 + " Sincerely," + [char]11 + " [Company Name] Customer Support"
